$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expenses")

# New expense entry: Bill_31 - Courier Load Box, Rs 500, paid by Akshay
$ws.Range("A33").Value = 33
$ws.Range("B33").Value = "8/12/2021"
$ws.Range("C33").Value = "Courier Load Box"
$ws.Range("D33").Value = 31
$ws.Range("E33").Value = 500
$ws.Range("F33").Value = "Akshay"
$ws.Range("G33").Value = "Bill_31"

# Link the bill reference to its scanned copy, mirroring the other bill rows
$ws.Hyperlinks.Add($ws.Range("G33"), "https://github.com/Akshay1595/Lithium-ion_battery_as_service/blob/master/Docs/Bills/Bill_31.jpeg")

# Re-apply the standard "Bill Link" cell format (hyperlink creation swaps in
# the built-in Hyperlink style by default) so G33 matches the rest of the column
$ws.Range("G32").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the cursor position left by the author on save
$ws.Range("N31").Select()
